$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(
    @(-18.93942138795675,1.811440133621348,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675),
    @(-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,2.617305898466031,-18.93942138795675,-18.93942138795675),
    @(-18.93942138795675,2.033082383869957,1.920787350772417,-18.93942138795675,4.321925368644031,-18.93942138795675,1.703977723360735,-18.93942138795675,2.123421226866288,-18.93942138795675),
    @(-18.93942138795675,1.777974737437811,-18.93942138795675,-18.93942138795675,-18.93942138795675,2.80686206474048,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675),
    @(-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675),
    @(2.645363615750298,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675),
    @(-18.93942138795675,-18.93942138795675,-18.93942138795675,1.996714010248535,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675),
    @(3.780670240346972,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675),
    @(-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,1.080561548030158,-18.93942138795675,2.028742647716157),
    @(-18.93942138795675,-18.93942138795675,-18.93942138795675,2.934785917340463,-18.93942138795675,2.9838571267412,-18.93942138795675,-18.93942138795675,-18.93942138795675,1.999084984374655),
    @(-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675),
    @(-18.93942138795675,-18.93942138795675,-18.93942138795675,2.449048455987207,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,1.863124164231722,1.907997694065814),
    @(-18.93942138795675,-18.93942138795675,0.9900985432516622,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,2.11358770015759),
    @(-18.93942138795675,-18.93942138795675,0.6760879188363068,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675),
    @(-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,1.956332119914588,-18.93942138795675),
    @(-18.93942138795675,2.584516189127088,2.938137230462033,-18.93942138795675,-18.93942138795675,-18.93942138795675,1.488150290702825,2.219383656914345,2.268584887367918,-18.93942138795675),
    @(-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,-18.93942138795675,1.848727386000359,1.373874647287896,1.724877448534649,-18.93942138795675),
    @(-18.93942138795675,-18.93942138795675,1.384135907034331,-18.93942138795675,-18.93942138795675,-18.93942138795675,1.447860778466786,1.460616850711499,-18.93942138795675,-18.93942138795675),
    @(-18.93942138795675,0.5029086280375464,1.235841850448465,-18.93942138795675,-18.93942138795675,-18.93942138795675,1.815114666259123,0.8058869249635578,-18.93942138795675,1.941660450451589),
    @(-18.93942138795675,0.6346935020825177,-18.93942138795675,1.537088499949748,-18.93942138795675,2.348054843664391,2.031384559713528,-18.93942138795675,-18.93942138795675,-18.93942138795675)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $rowValues = $values[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
